# "added comments and error page"
# Renames the three address-component headers to their snake_case form,
# corrects the apartment number for the McKinley Dr listing, restyles the
# header row (bold Calibri 11, thin box border, centered/top aligned),
# drops the header row's explicit 18pt height back to the sheet default,
# and leaves the cursor parked on F18 (as it was after the edit was made).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: "street number"/"street name"/"apt number" -> snake_case
$ws.Range("A1").Value = "street_number"
$ws.Range("B1").Value = "street_name"
$ws.Range("C1").Value = "apt_number"

# --- Data correction: McKinley Dr apartment number 1210 -> 1201
$ws.Range("C4").Value = 1201

# --- Header row styling (A1:F1): bold 11pt Calibri, boxed, centered + top
$header = $ws.Range("A1:F1")
$header.Style = "Normal"
$header.Font.Name = "Calibri"
$header.Font.Size = 11
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.LineStyle = 1         # xlContinuous
$header.Borders.Weight = 2            # xlThin

# --- Row 1 no longer carries a custom 18pt height
$ws.Rows.Item(1).AutoFit()

# --- Leave the selection where the editor left it
[void]$ws.Range("F18").Select()
